$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.069.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.678.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("E7").Value = "  +1.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07364"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.110"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.675.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.681"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001059"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06563"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9986"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.938"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.041.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.444"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.408"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.861.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.28%  "

$ws.Range("E31").Value = "  +2.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.097"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.849"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08494"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.674"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.195"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06119"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.236"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.35%  "

$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9986"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.837"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.973"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07033"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("E51").Value = "  +3.63%  "
